$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.108.49'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.424.79'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.531'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = '2.423.58'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('E12').Value = '  -3.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.350'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000174'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.80%  '
$ws.Range('D16').Value = '2.859.25'
$ws.Range('D17').Value = '61.984.22'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '2.409.31'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '323.26'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.83'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.43'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.60'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '558.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.58%  '
$ws.Range('D28').Value = '2.543.25'
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '0.0₃0936'
$ws.Range('E30').Value = '  -0.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.21'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.57%  '
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.76'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.50'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '151.31'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.68'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('E44').Value = '  -3.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.65'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0530'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0229'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.24%  '
